$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Selection Type value in B3
$ws.Range("B3").Value = "TRChooseByHighestDemand"

# Row 8: set A8:Y8
$row8 = New-Object 'object[,]' 1,25
$row8[0,0] = -1.0
$row8[0,1] = -1.0
$row8[0,2] = 0.0
$row8[0,3] = 108.0
$row8[0,4] = -1.0
$row8[0,5] = 37.0
$row8[0,6] = 38.0
$row8[0,7] = 36.0
$row8[0,8] = 34.0
$row8[0,9] = 31.0
$row8[0,10] = 30.0
$row8[0,11] = 18.0
$row8[0,12] = 25.0
$row8[0,13] = 26.0
$row8[0,14] = 27.0
$row8[0,15] = 29.0
$row8[0,16] = 28.0
$row8[0,17] = 57.0
$row8[0,18] = 58.0
$row8[0,19] = 49.0
$row8[0,20] = 44.0
$row8[0,21] = 46.0
$row8[0,22] = 39.0
$row8[0,23] = 40.0
$row8[0,24] = -1.0
$ws.Range("A8:Y8").Value = $row8
$ws.Range("Z8:Z8").ClearContents()

# Row 9: set A9:Y9
$row9 = New-Object 'object[,]' 1,25
$row9[0,0] = -1.0
$row9[0,1] = -1.0
$row9[0,2] = 1.0
$row9[0,3] = 108.0
$row9[0,4] = -1.0
$row9[0,5] = 36.0
$row9[0,6] = 34.0
$row9[0,7] = 33.0
$row9[0,8] = 31.0
$row9[0,9] = 30.0
$row9[0,10] = 18.0
$row9[0,11] = 19.0
$row9[0,12] = 22.0
$row9[0,13] = 26.0
$row9[0,14] = 27.0
$row9[0,15] = 29.0
$row9[0,16] = 28.0
$row9[0,17] = 53.0
$row9[0,18] = 50.0
$row9[0,19] = 57.0
$row9[0,20] = 49.0
$row9[0,21] = 44.0
$row9[0,22] = 39.0
$row9[0,23] = 40.0
$row9[0,24] = -1.0
$ws.Range("A9:Y9").Value = $row9

# Row 10: set A10:X10
$row10 = New-Object 'object[,]' 1,24
$row10[0,0] = -1.0
$row10[0,1] = -1.0
$row10[0,2] = 2.0
$row10[0,3] = 108.0
$row10[0,4] = -1.0
$row10[0,5] = 36.0
$row10[0,6] = 31.0
$row10[0,7] = 30.0
$row10[0,8] = 11.0
$row10[0,9] = 18.0
$row10[0,10] = 25.0
$row10[0,11] = 26.0
$row10[0,12] = 27.0
$row10[0,13] = 29.0
$row10[0,14] = 28.0
$row10[0,15] = 38.0
$row10[0,16] = 56.0
$row10[0,17] = 57.0
$row10[0,18] = 50.0
$row10[0,19] = 49.0
$row10[0,20] = 44.0
$row10[0,21] = 46.0
$row10[0,22] = 5.0
$row10[0,23] = -1.0
$ws.Range("A10:X10").Value = $row10
$ws.Range("Y10:AA10").ClearContents()

# Row 11: set A11:Y11
$row11 = New-Object 'object[,]' 1,25
$row11[0,0] = -1.0
$row11[0,1] = -1.0
$row11[0,2] = 3.0
$row11[0,3] = 108.0
$row11[0,4] = -1.0
$row11[0,5] = 35.0
$row11[0,6] = 36.0
$row11[0,7] = 34.0
$row11[0,8] = 30.0
$row11[0,9] = 18.0
$row11[0,10] = 22.0
$row11[0,11] = 26.0
$row11[0,12] = 27.0
$row11[0,13] = 29.0
$row11[0,14] = 28.0
$row11[0,15] = 52.0
$row11[0,16] = 51.0
$row11[0,17] = 50.0
$row11[0,18] = 57.0
$row11[0,19] = 49.0
$row11[0,20] = 39.0
$row11[0,21] = 40.0
$row11[0,22] = 6.0
$row11[0,23] = 4.0
$row11[0,24] = -1.0
$ws.Range("A11:Y11").Value = $row11

# Row 12: set A12:X12
$row12 = New-Object 'object[,]' 1,24
$row12[0,0] = -1.0
$row12[0,1] = -1.0
$row12[0,2] = 4.0
$row12[0,3] = 108.0
$row12[0,4] = -1.0
$row12[0,5] = 37.0
$row12[0,6] = 38.0
$row12[0,7] = 36.0
$row12[0,8] = 31.0
$row12[0,9] = 30.0
$row12[0,10] = 33.0
$row12[0,11] = 18.0
$row12[0,12] = 25.0
$row12[0,13] = 26.0
$row12[0,14] = 27.0
$row12[0,15] = 29.0
$row12[0,16] = 28.0
$row12[0,17] = 56.0
$row12[0,18] = 57.0
$row12[0,19] = 58.0
$row12[0,20] = 50.0
$row12[0,21] = 49.0
$row12[0,22] = 44.0
$row12[0,23] = -1.0
$ws.Range("A12:X12").Value = $row12
$ws.Range("Y12:Z12").ClearContents()

# Row 13: set A13:Y13
$row13 = New-Object 'object[,]' 1,25
$row13[0,0] = -1.0
$row13[0,1] = -1.0
$row13[0,2] = 5.0
$row13[0,3] = 108.0
$row13[0,4] = -1.0
$row13[0,5] = 35.0
$row13[0,6] = 36.0
$row13[0,7] = 34.0
$row13[0,8] = 33.0
$row13[0,9] = 30.0
$row13[0,10] = 18.0
$row13[0,11] = 22.0
$row13[0,12] = 26.0
$row13[0,13] = 27.0
$row13[0,14] = 29.0
$row13[0,15] = 28.0
$row13[0,16] = 53.0
$row13[0,17] = 52.0
$row13[0,18] = 51.0
$row13[0,19] = 50.0
$row13[0,20] = 57.0
$row13[0,21] = 49.0
$row13[0,22] = 7.0
$row13[0,23] = 40.0
$row13[0,24] = -1.0
$ws.Range("A13:Y13").Value = $row13
